# Auto-generated Excel COM-interop edit script
# Applies targeted cell value updates to multiple sheets per the source diff.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1378.2222
$ws.Range("I62").Value = 1450.5
$ws.Range("J62").Value = 800
$ws.Range("K62").Value = 1450.5
$ws.Range("L62").Value = 800
$ws.Range("M62").Value = -826.5
$ws.Range("N62").Value = -2048
$ws.Range("H65").Value = 1378.2222
$ws.Range("I65").Value = 1450.5
$ws.Range("J65").Value = 800
$ws.Range("K65").Value = 7252.5
$ws.Range("L65").Value = 4000
$ws.Range("M65").Value = -4132.5
$ws.Range("N65").Value = -10240
$ws.Range("H125").Value = 592.8333
$ws.Range("I125").Value = 559.1429000000001
$ws.Range("J125").Value = 640
$ws.Range("K125").Value = 5032.2861
$ws.Range("L125").Value = 5760
$ws.Range("M125").Value = -2572.2861
$ws.Range("N125").Value = -10680

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 14655.611
$ws.Range("I86").Value = 22135.818
$ws.Range("J86").Value = 2901
$ws.Range("K86").Value = 22135.818
$ws.Range("L86").Value = 2901
$ws.Range("M86").Value = -21012.818
$ws.Range("N86").Value = -5147
$ws.Range("H89").Value = 14655.611
$ws.Range("I89").Value = 22135.818
$ws.Range("J89").Value = 2901
$ws.Range("K89").Value = 110679.09
$ws.Range("L89").Value = 14505
$ws.Range("M89").Value = -105063.09
$ws.Range("N89").Value = -25737

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 35718696
$ws.Range("I58").Value = 50004068
$ws.Range("J58").Value = 5262.625
$ws.Range("K58").Value = 50004068
$ws.Range("L58").Value = 5262.625
$ws.Range("M58").Value = -50003865
$ws.Range("N58").Value = -5668.625
$ws.Range("H62").Value = 2821.5386
$ws.Range("I62").Value = 2821.5386
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2821.5386
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2197.5386
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 2821.5386
$ws.Range("I65").Value = 2821.5386
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 14107.693
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -10987.693
$ws.Range("N65").ClearContents()
$ws.Range("H136").Value = 35718696
$ws.Range("I136").Value = 50004068
$ws.Range("J136").Value = 5262.625
$ws.Range("K136").Value = 150012204
$ws.Range("L136").Value = 15787.875
$ws.Range("M136").Value = -150009654
$ws.Range("N136").Value = -20887.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 616.6667
$ws.Range("J39").Value = 616.6667
$ws.Range("L39").Value = 1850.0001
$ws.Range("N39").Value = -2438.0001
$ws.Range("H43").Value = 10000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 30000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -30228
$ws.Range("H44").Value = 678.7
$ws.Range("I44").Value = 399.5
$ws.Range("J44").Value = 748.5
$ws.Range("K44").Value = 1198.5
$ws.Range("L44").Value = 2245.5
$ws.Range("M44").Value = -800.5
$ws.Range("N44").Value = -3041.5
$ws.Range("H46").Value = 3217.7778
$ws.Range("J46").Value = 3851.4285
$ws.Range("L46").Value = 11554.2855
$ws.Range("N46").Value = -11736.2855
$ws.Range("H47").Value = 1072
$ws.Range("I47").Value = 595
$ws.Range("J47").Value = 2980
$ws.Range("K47").Value = 1785
$ws.Range("L47").Value = 8940
$ws.Range("M47").Value = -1354
$ws.Range("N47").Value = -9802
$ws.Range("H48").Value = 5000
$ws.Range("J48").Value = 5000
$ws.Range("L48").Value = 15000
$ws.Range("N48").Value = -15500
$ws.Range("H49").Value = 3188
$ws.Range("I49").Value = 500
$ws.Range("J49").Value = 3636
$ws.Range("K49").Value = 1500
$ws.Range("L49").Value = 10908
$ws.Range("M49").Value = -1344
$ws.Range("N49").Value = -11220
$ws.Range("H54").Value = 4000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H55").Value = 3428.5715
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3428.5715
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 10285.7145
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -10639.7145
$ws.Range("H140").Value = 2048.39
$ws.Range("I140").Value = 1169.8125
$ws.Range("J140").Value = 2215.738
$ws.Range("K140").Value = 3509.4375
$ws.Range("L140").Value = 6647.214
$ws.Range("M140").Value = 1670.5625
$ws.Range("N140").Value = -17007.214

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1470.35
$ws.Range("I102").Value = 1462.6875
$ws.Range("J102").Value = 1501
$ws.Range("K102").Value = 1462.6875
$ws.Range("L102").Value = 1501
$ws.Range("M102").Value = 159.3125
$ws.Range("N102").Value = -4745
$ws.Range("H122").Value = 1383.3334
$ws.Range("I122").Value = 1600
$ws.Range("K122").Value = 4800
$ws.Range("M122").Value = -2350
$ws.Range("H126").Value = 1960
$ws.Range("I126").Value = 1900
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 5700
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3230
$ws.Range("N126").Value = -10940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2599.2222
$ws.Range("I7").Value = 2012
$ws.Range("J7").Value = 3333.25
$ws.Range("K7").Value = 2012
$ws.Range("L7").Value = 3333.25
$ws.Range("M7").Value = -1900
$ws.Range("N7").Value = -3557.25
$ws.Range("H40").Value = 2457.1428
$ws.Range("I40").Value = 2640
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 2640
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -2504
$ws.Range("N40").Value = -2272
$ws.Range("H122").Value = 4320.737
$ws.Range("I122").Value = 3944
$ws.Range("J122").Value = 4594.727
$ws.Range("K122").Value = 11832
$ws.Range("L122").Value = 13784.181
$ws.Range("M122").Value = -9382
$ws.Range("N122").Value = -18684.181
$ws.Range("H126").Value = 2599.2222
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 3333.25
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 9999.75
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -14939.75
$ws.Range("H132").Value = 49297.727
$ws.Range("I132").Value = 2955.9
$ws.Range("J132").Value = 87915.914
$ws.Range("K132").Value = 8867.700000000001
$ws.Range("L132").Value = 263747.742
$ws.Range("M132").Value = -6337.700000000001
$ws.Range("N132").Value = -268807.742

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2028
$ws.Range("I126").Value = 1186.8572
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 3560.5716
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -1090.5716
$ws.Range("N126").Value = -15440

Write-Host "Applied profit/price updates across ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets."